$wb = $excel.ActiveWorkbook
$ws2 = $wb.Worksheets.Item("2018-F")
$ws2.Activate()
Write-Host "Before ScrollRow:" $excel.ActiveWindow.ScrollRow
$excel.ActiveWindow.ScrollRow = 4
Write-Host "After ScrollRow:" $excel.ActiveWindow.ScrollRow
